$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Credentials")

# Update credentials in the "Login As Host" block (B2:B4)
$ws.Range("B2").Value = "https://stagingtshq.bsbtest.com/Default.aspx?portalid=24904"
$ws.Range("B3").Value = "vipul24904"
$ws.Range("B4").Value = "St4ckSp0rts@"

# Match styles used by the other "Login As Host"-style entries (B15:B24 use style referenced as s=19/20)
$ws.Range("B2").Style = $ws.Range("B16").Style
$ws.Range("B3").Style = $ws.Range("B16").Style
$ws.Range("B4").Style = $ws.Range("B25").Style

# Update the password for the bottom admin credential block
$ws.Range("B25").Value = "#Sp0rtsC0nn3ct"

# Update the active selection to B4 like in the updated file
$ws.Range("B4").Select()

# Remove file-recovery marker left over from a crash/repair
$wb.Save()
